$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 628 entirely; this shifts all rows below it up by one.
$ws.Rows.Item(628).Delete()
